# Sprint 4 backlog: add "Reviewers 1" (col E) / "Reviewers 2" (col F) assignments
# for each backlog row, and move the saved selection to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: row, Reviewers 1, Reviewers 2
$reviewerRows = @(
    @(2,  "fatayry",           "reyad"),
    @(3,  "bassem",            "khaled"),
    @(4,  "youssef abo bakr",  "omar fatayry"),
    @(5,  "bassem",            "ahmed hesham"),
    @(6,  "youssef waly",      "ahmed hesham"),
    @(7,  "youssef abo bakr",  "atef"),
    @(8,  "reyad",             "andrew"),
    @(9,  "andrew",            "khaled"),
    @(10, "atef",              "alaa"),
    @(11, "alaa",              "youssef waly")
)

foreach ($entry in $reviewerRows) {
    $row = $entry[0]
    $ws.Cells.Item($row, 5).Value = $entry[1]
    $ws.Cells.Item($row, 6).Value = $entry[2]
}

# Move the active cell / selection to E6, matching the saved view state.
$ws.Range("E6").Select() | Out-Null
